# HarpLakeparametres.xlsx - "Formatting edits for Harp input files"
# The parameter table on the first sheet (Hoja1) had a blank row above the
# header row. Delete that leading blank row so the header ("Parameter",
# "Value", "Units") starts at row 1 and everything below shifts up by one,
# closing the former row-19 gap up into row 18 (now row 17->18 gap moves to
# 18->19 vacancy as before, i.e. the blank separator row keeps its relative
# position one row higher).

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Hoja1")

$ws.Rows("1:1").Delete()

# Leave the selection where the author left it after the edit.
$ws.Range("E4").Select()
